# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.136.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.602.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.82%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.257"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.58%  "

$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08140"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.590"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001259"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.397"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.603.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06863"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.558"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  -5.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.148.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.347"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.720"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.304"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.424"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.848"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -13.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.780.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07688"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9492"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02737"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.262"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2549"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08907"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.95%  "

$ws.Range("E40").Value = "  -6.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.385"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7107"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6635"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("E47").Value = "  -5.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.980"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.208"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
